# Update workbook to reflect data through 2022-08-08 (commit: "Add data for 2022-08-16")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (from "Through 2022-08-07" to "Through 2022-08-08")
$ws.Name = "Through 2022-08-08"

# Update the August row label in column A (row 9)
$ws.Range("A9").Value = "August (through 08-08)"

# Update August row (row 9) values for columns D..I (2017, 2018, 2019, 2020, 2021, 2022)
$ws.Range("D9").Value = 22
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = 13
$ws.Range("G9").Value = 44
$ws.Range("H9").Value = 53
$ws.Range("I9").Value = 42

# Update Total row (row 10) values for columns D..I
$ws.Range("D10").Value = 487
$ws.Range("E10").Value = 441
$ws.Range("F10").Value = 317
$ws.Range("G10").Value = 665
$ws.Range("H10").Value = 963
$ws.Range("I10").Value = 1012

$wb.Save()
